$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1825842696629214
$ws.Range("C2").Value = 0.5589887640449438
$ws.Range("J2").Value = 0.02247191011235955
$ws.Range("P2").Value = 0.1292134831460674
$ws.Range("S2").Value = 0.1067415730337079
$ws.Range("C3").Value = 0.03349282296650718
$ws.Range("J3").Value = 0.01913875598086124
$ws.Range("P3").Value = 0.6985645933014354
$ws.Range("S3").Value = 0.2488038277511962
$ws.Range("P4").Value = 0.6
$ws.Range("S4").Value = 0.4
$ws.Range("B6").Value = 0.05019305019305019
$ws.Range("D6").Value = 0.01158301158301158
$ws.Range("F6").Value = 0.05405405405405406
$ws.Range("J6").Value = 0.2239382239382239
$ws.Range("O6").Value = 0.0193050193050193
$ws.Range("Q6").Value = 0.1698841698841699
$ws.Range("R6").Value = 0.09652509652509653
$ws.Range("S6").Value = 0.3745173745173745
$ws.Range("B7").Value = 0.1350210970464135
$ws.Range("D7").Value = 0.02109704641350211
$ws.Range("F7").Value = 0.05907172995780591
$ws.Range("J7").Value = 0.09282700421940929
$ws.Range("O7").Value = 0.01265822784810127
$ws.Range("Q7").Value = 0.1856540084388186
$ws.Range("R7").Value = 0.07172995780590717
$ws.Range("S7").Value = 0.4219409282700422
$ws.Range("B8").Value = 0.08823529411764706
$ws.Range("D8").Value = 0.009803921568627451
$ws.Range("F8").Value = 0.09215686274509804
$ws.Range("J8").Value = 0.1098039215686274
$ws.Range("O8").Value = 0.03529411764705882
$ws.Range("Q8").Value = 0.1745098039215686
$ws.Range("R8").Value = 0.09019607843137255
$ws.Range("S8").Value = 0.4
$ws.Range("B9").Value = 0.07722007722007722
$ws.Range("D9").Value = 0.007722007722007722
$ws.Range("E9").Value = 0.003861003861003861
$ws.Range("F9").Value = 0.05791505791505792
$ws.Range("J9").Value = 0.1351351351351351
$ws.Range("O9").Value = 0.01544401544401544
$ws.Range("Q9").Value = 0.1158301158301158
$ws.Range("R9").Value = 0.1235521235521236
$ws.Range("S9").Value = 0.4633204633204633
$ws.Range("B10").Value = 0.1088807785888078
$ws.Range("D10").Value = 0.01642335766423358
$ws.Range("E10").Value = 0.0006082725060827251
$ws.Range("F10").Value = 0.07116788321167883
$ws.Range("J10").Value = 0.1192214111922141
$ws.Range("O10").Value = 0.0218978102189781
$ws.Range("Q10").Value = 0.2141119221411192
$ws.Range("R10").Value = 0.09610705596107055
$ws.Range("S10").Value = 0.3515815085158151
$ws.Range("G11").Value = 0.1343669250645995
$ws.Range("J11").Value = 0.1085271317829457
$ws.Range("K11").Value = 0.2144702842377261
$ws.Range("L11").Value = 0.5348837209302325
$ws.Range("S11").Value = 0.007751937984496124
$ws.Range("G12").Value = 0.7129186602870813
$ws.Range("J12").Value = 0.1961722488038277
$ws.Range("K12").Value = 0.004784688995215311
$ws.Range("L12").Value = 0.02392344497607655
$ws.Range("S12").Value = 0.06220095693779904
$ws.Range("F13").Value = 0.01818181818181818
$ws.Range("G13").Value = 0.7272727272727273
$ws.Range("J13").Value = 0.2181818181818182
$ws.Range("S13").Value = 0.03636363636363636
$ws.Range("G14").Value = 0.5
$ws.Range("J14").Value = 0.5
$ws.Range("F15").Value = 0.01550387596899225
$ws.Range("H15").Value = 0.1279069767441861
$ws.Range("I15").Value = 0.09302325581395349
$ws.Range("J15").Value = 0.3643410852713178
$ws.Range("K15").Value = 0.03875968992248062
$ws.Range("M15").Value = 0.003875968992248062
$ws.Range("O15").Value = 0.05038759689922481
$ws.Range("S15").Value = 0.3062015503875969
$ws.Range("F16").Value = 0.009523809523809525
$ws.Range("H16").Value = 0.119047619047619
$ws.Range("I16").Value = 0.1095238095238095
$ws.Range("J16").Value = 0.4666666666666667
$ws.Range("K16").Value = 0.1047619047619048
$ws.Range("M16").Value = 0.01904761904761905
$ws.Range("O16").Value = 0.02857142857142857
$ws.Range("S16").Value = 0.1428571428571428
$ws.Range("F17").Value = 0.007220216606498195
$ws.Range("H17").Value = 0.1967509025270758
$ws.Range("I17").Value = 0.09386281588447654
$ws.Range("J17").Value = 0.4205776173285198
$ws.Range("K17").Value = 0.09025270758122744
$ws.Range("M17").Value = 0.01985559566787004
$ws.Range("N17").Value = 0.001805054151624549
$ws.Range("O17").Value = 0.0631768953068592
$ws.Range("S17").Value = 0.1064981949458484
$ws.Range("F18").Value = 0.003597122302158274
$ws.Range("H18").Value = 0.1798561151079137
$ws.Range("I18").Value = 0.1223021582733813
$ws.Range("J18").Value = 0.4568345323741007
$ws.Range("K18").Value = 0.06115107913669065
$ws.Range("M18").Value = 0.007194244604316547
$ws.Range("O18").Value = 0.04316546762589928
$ws.Range("S18").Value = 0.1258992805755396
$ws.Range("F19").Value = 0.007518796992481203
$ws.Range("H19").Value = 0.1860902255639098
$ws.Range("I19").Value = 0.08270676691729323
$ws.Range("J19").Value = 0.3984962406015037
$ws.Range("K19").Value = 0.1234335839598997
$ws.Range("M19").Value = 0.02380952380952381
$ws.Range("N19").Value = 0.001879699248120301
$ws.Range("O19").Value = 0.05889724310776942
$ws.Range("S19").Value = 0.1171679197994987
